$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("I10").Value = 500.4
$ws1.Range("M10").Value = 8649.1
$ws1.Range("M12").Value = 3399.19
$ws1.Range("M18").Value = 1702.17
$ws1.Range("M22").Value = "9 de 20"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = 10041.73
$ws2.Range("F12").Value = 3399.19
$ws2.Range("F18").Value = 1702.17
$ws2.Range("F22").Value = 49964.29

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D8").Value = 1188
$ws3.Range("E8").Value = -563
$ws3.Range("F8").Value = 1.9008

$ws3.Range("D16").Value = 40490.24
$ws3.Range("E16").Value = 3776
$ws3.Range("F16").Value = 0.9146979729925108

$ws3.Range("D19").Value = 49964.28999999999
$ws3.Range("E19").Value = 15413.70762291768
$ws3.Range("F19").Value = 0.7642370800063392
